# Regenerate the "K" column (G) of the save_data sheet for leiter_mark.
# The K column previously held a "Strike#" style count; this recalculates
# it (e.g. strikeouts derived from the underlying pitch-by-pitch values)
# and overwrites each row's G cell with the freshly computed value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (recalculated s_val for that game).
$kValues = [ordered]@{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 3
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 4
    25 = 1
    26 = 2
    27 = 0
    28 = 1
    29 = 1
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 1
    35 = 2
    36 = 0
    37 = 2
    39 = 1
    40 = 1
    41 = 1
    42 = 2
    43 = 1
    44 = 2
    45 = 2
    46 = 2
    47 = 3
    48 = 1
    49 = 0
    50 = 1
    51 = 1
    52 = 2
    53 = 3
    54 = 0
    55 = 2
    56 = 1
    57 = 2
    58 = 3
    59 = 1
    60 = 0
    61 = 1
    62 = 1
    63 = 2
    64 = 2
    65 = 3
    66 = 2
    67 = 0
    68 = 1
    69 = 2
    70 = 2
    71 = 1
    72 = 2
    73 = 3
    74 = 2
    76 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
